$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Python code snippets (matching existing pythonCode column style: wrap text, s=1)
$groceryAppend = "grocery = [`"apples`", `"bananas`", `"cucumbers`", `"dates`", `"strawberries`"]`ngrocery.append(`"oranges`")`nfor fruits in grocery:`n  print(fruits)"
$groceryInsertPop = "grocery_list = [`"apples`", `"bananas`", `"cucumbers`", `"dates`", `"strawberries`"]`ngrocery_list.insert(2,`"cashews`")`ngrocery_list.pop()`nfor list in grocery_list:`n  print(list)"
$marksAverage = "marks = [3, 5, 4, 2, 5, 5, 3, 5, 4, 4, 4]`nsum_of_marks=sum(marks)`nlen_of_marks=len(marks)`nprint(sum_of_marks)`nprint(len_of_marks)`nresult = sum(marks) / len(marks)`nprint(`"The final grade:`" + str(result))"

$ws.Range("A4").Value = $groceryAppend
$ws.Range("A5").Value = $groceryInsertPop
$ws.Range("A6").Value = $marksAverage

# Match formatting of the other python-code cells (wrap text style)
$ws.Range("A4:A6").WrapText = $true

$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 75
$ws.Rows.Item(6).RowHeight = 105

# Update selection as in the target sheet view
$ws.Range("A8:A14").Select() | Out-Null
